# Curated ICDC 1-8, 29, 30, 31: update startup sheet queries, labels, and formatting
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (unchanged text) ----
$ws.Range("A1").Value = "TabName"
$ws.Range("B1").Value = "query"
$ws.Range("C1").Value = "StatQuery"
$ws.Range("D1").Value = "dbExcel"
$ws.Range("E1").Value = "WebExcel"

# ---- Column A tab names ----
$ws.Range("A2").Value = "CasesTab"
$ws.Range("A3").Value = "SamplesTab"
$ws.Range("A4").Value = "FilesTab"
$ws.Range("A5").Value = "StudyFilesTab"

# ---- Column D / E output file names ----
$neo4jFile = "TC04_Canine_Filter_Study-UBC01_Neo4jData.xlsx"
$webFile = "TC04_Canine_Filter_Study-UBC01_WebData.xlsx"
$ws.Range("D2").Value = $neo4jFile
$ws.Range("E2").Value = $webFile
$ws.Range("D3").Value = $neo4jFile
$ws.Range("E3").Value = $webFile
$ws.Range("D4").Value = $neo4jFile
$ws.Range("E4").Value = $webFile
$ws.Range("D5").Value = $neo4jFile
$ws.Range("E5").Value = $webFile

# ---- Column B: per-tab Cypher queries ----
$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE s.clinical_study_designation IN ['UBC01']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co, demo.patient_age_at_enrollment AS age
RETURN  
       coalesce(c.case_id, '') AS `Case ID`,
       coalesce(s.clinical_study_designation, '') AS `Study Code`,
       coalesce(s.clinical_study_type, '') AS  `Study Type`,
       coalesce(demo.breed, '') AS Breed ,
       coalesce(diag.disease_term, '') AS Diagnosis ,
       coalesce(diag.stage_of_disease, '') AS `Stage of Disease`,
     coalesce(CASE age % 1 WHEN 0 THEN apoc.convert.toInteger(age) ELSE age END, '') AS Age,
       coalesce(demo.sex, '') AS Sex,
       coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
coalesce(CASE weight % 1 WHEN 0 THEN apoc.convert.toInteger(weight) ELSE weight END, '') AS `Weight (kg)`,
       coalesce(diag.best_response, '') AS `Response to Treatment`,
       coalesce(co.cohort_description, '') AS `Cohort`
'@
$ws.Range("B2").Value = $casesQuery

$samplesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) 
WHERE s.clinical_study_designation IN ['UBC01']
WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '') AS `Sample ID`, 
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(demo.breed,'') AS Breed , 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(samp.sample_site, '') AS `Sample Site`,
        coalesce(samp.summarized_sample_type, '') AS `Sample Type`,
        coalesce(samp.specific_sample_pathology, '') AS `Pathology/Morphology`,
        coalesce(samp.tumor_grade, '') AS `Tumor Grade`,
        coalesce(samp.sample_chronology, '') AS `Sample Chronology`,
        coalesce(samp.percentage_tumor, '') AS `Percentage Tumor`,
        coalesce(samp.necropsy_sample, '') AS `Necropsy Sample`,
        coalesce(samp.sample_preservation, '') AS `Sample Preservation`
'@
$ws.Range("B3").Value = $samplesQuery

$filesQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f)-[*]->(samp:sample)-->(c)-->(s:study)
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE s.clinical_study_designation IN ['UBC01']
OPTIONAL MATCH (s:study)<--(c)<--(diag:diagnosis)<-[*]-(samp)
WITH
        f, parent, c, demo, diag, s, samp,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent, c, demo, diag, s, samp,
        f.file_size /(1024^i) AS value, 
        10^precision AS factor,
        units[i] as unit
WITH    
        f, parent, c, demo, diag, s, samp, unit,
        round(factor * value)/factor AS size
RETURN 
        coalesce(f.file_name, '') AS `File Name`,
        coalesce(f.file_type, '') AS `File Type`,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
        coalesce(samp.sample_id, '') AS `Sample ID`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(demo.breed,'') AS Breed ,
        coalesce(diag.disease_term,'') AS Diagnosis
'@
$ws.Range("B4").Value = $filesQuery

$studyFilesQuery = @'
MATCH (f:file)-->(s:study)
MATCH (s)<--(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)
WHERE s.clinical_study_designation IN ['UBC01']
WITH DISTINCT f,  s, c, demo, diag
WITH
        f, c, demo, diag, s,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH    
        f, c, demo, diag, s,
        f.file_size /(1024^i) AS value, 10^precision AS factor,
        units[i] as unit
        WITH    
        f,  c, demo, diag, s, unit,
        round(factor * value)/factor AS size
RETURN DISTINCT
  coalesce(f.file_name, '') AS `File Name`,
  coalesce(f.file_type, '') AS `File Type`,
  coalesce("study", '') AS `Association`,
  coalesce(f.file_description, '') AS `Description`,
  coalesce(f.file_format, '') AS  Format,
  CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
  coalesce(s.clinical_study_designation,'') AS `Study Code`
'@
$ws.Range("B5").Value = $studyFilesQuery

# ---- Column C: shared StatQuery (same text on every data row) ----
$statQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE s.clinical_study_designation IN ['UBC01']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@
$ws.Range("C2").Value = $statQuery
$ws.Range("C3").Value = $statQuery
$ws.Range("C4").Value = $statQuery
$ws.Range("C5").Value = $statQuery

# ---- Formatting: font size 14 across used range, wrap text for query columns ----
$ws.Range("A1:E1").Font.Size = 14
$ws.Range("A2:A5").Font.Size = 14
$ws.Range("D2:E5").Font.Size = 14
$ws.Range("B2:C5").Font.Size = 14
$ws.Range("B2:C5").WrapText = $true
$ws.Range("C6").Font.Size = 14
$ws.Range("C6").WrapText = $true

# ---- Row heights ----
$ws.Rows("1").RowHeight = 38.25
$ws.Rows("2:6").RowHeight = 63

# ---- Column widths ----
$ws.Columns("A").ColumnWidth = 10.022135416666666
$ws.Columns("B").ColumnWidth = 57.166666666666664
$ws.Columns("C").ColumnWidth = 64.30729166666667
$ws.Columns("D").ColumnWidth = 54.592447916666664
$ws.Columns("E").ColumnWidth = 38.877604166666664

# ---- Selection ----
[void]$ws.Range("D2").Select()
